# Insert a new weekly price record above row 1119 (Hortaliza / Cebolla
# subset for "Terminal La Palmera de La Serena"), shifting the existing
# rows 1119:1175 down to 1120:1176.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh row at position 1119; Excel pushes rows 1119-1175 down
# to 1120-1176 and extends the used range accordingly.
$ws.Rows.Item(1119).Insert()

# Populate the new row 1119 with the new record's data.
$ws.Cells.Item(1119, 1).Value = 8
$ws.Cells.Item(1119, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(1119, 3).Value = "Coquimbo"
$ws.Cells.Item(1119, 4).Value = 45147
$ws.Cells.Item(1119, 5).Value = 4
$ws.Cells.Item(1119, 6).Value = 100112004
$ws.Cells.Item(1119, 7).Value = "Cebolla"
$ws.Cells.Item(1119, 8).Value = "Sin especificar"
$ws.Cells.Item(1119, 9).Value = "1a (guarda)"
$ws.Cells.Item(1119, 10).Value = 2000
$ws.Cells.Item(1119, 11).Value = 10500
$ws.Cells.Item(1119, 12).Value = 11000
$ws.Cells.Item(1119, 13).Value = 10750
$ws.Cells.Item(1119, 14).Value = "`$/malla 16 kilos"
$ws.Cells.Item(1119, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(1119, 16).Value = 672
$ws.Cells.Item(1119, 17).Value = 16
$ws.Cells.Item(1119, 18).Value = "Hortaliza"
